$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new test-candidate data
$ws.Range("A2").Value = "VOxrJ584"
$ws.Range("B2").Value = 2012454477
$ws.Range("C2").Value = "kcmoyof93"
$ws.Range("D2").Value = "S&Q49kt$"
$ws.Range("E2").Value = "MR"
$ws.Range("F2").Value = "EveOCaMs"
$ws.Range("G2").Value = "rycI"
$ws.Range("H2").Value = "Candidate"

# Remove row 3 entirely (was the second candidate row)
$ws.Rows("3:3").Delete()

# Refresh the active selection to reflect the new used range
$ws.Range("A1:H2").Select() | Out-Null
